# "Generate Report for Handoff"
#
# This updates the localization-status report so it reflects a fresh
# handoff generation instead of the previous handback state:
#   - Status "Handed back: in sync with en-US" -> "Ready for handoff"
#     (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
#   - The "Latest HO Xliff Generate Date" on the Overview sheet and the
#     matching "Latest Handoff Datetime" on the de-de sheet move from
#     2016-08-18 12:59:11 -> 2016-08-18 12:59:58
#   - The "Latest Handoff Datetime" on the zh-cn sheet moves from
#     2016-08-18 12:59:00 -> 2016-08-18 12:59:52
#   - Columns E/F on the Overview sheet and column C on the zh-cn/de-de
#     sheets (the "Status" columns) are narrowed to fit the shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Datetimes
$overview.Range("G2").Value = "2016-08-18 12:59:58"
$dede.Range("H2").Value = "2016-08-18 12:59:58"
$zhcn.Range("H2").Value = "2016-08-18 12:59:52"

# --- Narrow the Status columns now that the text is shorter
$overview.Range("E1:F1").ColumnWidth = 16.333333333333332
$zhcn.Range("C1").ColumnWidth = 16.333333333333332
$dede.Range("C1").ColumnWidth = 16.333333333333332
